$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 91: the date/time stamp was corrected from an intraday timestamp
#     to the 07:00 (daily open) timestamp used by the rest of the sheet.
$ws.Range("A91").Value2 = 45474.2916666667

# --- Row 92: new trading-data row appended for 2024-07-02.
# Copy the date cell's formatting (style index) from the row above so the
# new date cell keeps the existing custom date style instead of creating a
# brand-new style entry.
$ws.Range("A91").Copy()
$ws.Range("A92").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("A92").Value2 = 45475.622037037

$ws.Range("B92").Value2 = 6000
$ws.Range("C92").Value2 = 3.25999999046326
$ws.Range("D92").Value2 = 3.20000004768372
$ws.Range("E92").Value2 = 3.23000001907349
$ws.Range("F92").Value2 = 3.20000004768372

# adj_close (column G) is stored as text/shared-string in this sheet, so
# force text formatting before assignment, then restore the default
# "Normal" style so no extra number-format is left applied to the cell.
$ws.Range("G92").NumberFormat = "@"
$ws.Range("G92").Formula = "3.20000004768372"
$ws.Range("G92").Style = "Normal"

# ticker (column H) is also text; copy formatting from the cell above then
# set the value.
$ws.Range("H91").Copy()
$ws.Range("H92").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("H92").Value2 = "ESPE.MI"

$excel.CutCopyMode = $false
